$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for team record columns, matching the style of the
# existing header row (bold, centered, bordered - same as A1:AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the team record (W/L/T) for every data row.
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 30).Value = 73
    $ws.Cells.Item($r, 31).Value = 89
    $ws.Cells.Item($r, 32).Value = 0
}
